# Auto-generated edit script: update crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Cell, [string]$Text)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (prevents strings that look numeric, e.g. '22.50' or '1.00', from
    # being silently parsed into numbers and losing their exact formatting).
    $Cell.Value = "'" + $Text
    # Reset the style back to Normal so the quote-prefix flag introduced
    # by the assignment above doesn't leave a stray style on the cell.
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '39.551.69'
Set-TextCell $ws.Range("E2") '  +2.01%  '
Set-TextCell $ws.Range("D3") '2.175.79'
Set-TextCell $ws.Range("E3") '  +3.82%  '
Set-TextCell $ws.Range("E4") '  +0.06%  '
Set-TextCell $ws.Range("D5") '229.96'
Set-TextCell $ws.Range("E5") '  +0.38%  '
Set-TextCell $ws.Range("E6") '  +1.35%  '
Set-TextCell $ws.Range("D7") '65.15'
Set-TextCell $ws.Range("E7") '  +6.56%  '
Set-TextCell $ws.Range("D9") '0.402'
Set-TextCell $ws.Range("E9") '  +4.13%  '
Set-TextCell $ws.Range("D10") '0.0867'
Set-TextCell $ws.Range("E10") '  +2.75%  '
Set-TextCell $ws.Range("E11") '  +0.26%  '
Set-TextCell $ws.Range("D12") '16.12'
Set-TextCell $ws.Range("E12") '  +5.56%  '
Set-TextCell $ws.Range("D13") '2.499.00'
Set-TextCell $ws.Range("E13") '  +3.94%  '
Set-TextCell $ws.Range("D14") '22.50'
Set-TextCell $ws.Range("E14") '  +2.16%  '
Set-TextCell $ws.Range("D15") '0.821'
Set-TextCell $ws.Range("E15") '  +0.29%  '
Set-TextCell $ws.Range("E16") '  +2.45%  '
Set-TextCell $ws.Range("D17") '2.168.22'
Set-TextCell $ws.Range("E17") '  +2.74%  '
Set-TextCell $ws.Range("D18") '39.505.34'
Set-TextCell $ws.Range("E18") '  +2.12%  '
Set-TextCell $ws.Range("D19") '72.59'
Set-TextCell $ws.Range("E19") '  +1.20%  '
Set-TextCell $ws.Range("D20") '6.19'
Set-TextCell $ws.Range("E20") '  +1.28%  '
Set-TextCell $ws.Range("D21") '0.0₃0859'
Set-TextCell $ws.Range("E21") '  +2.08%  '
Set-TextCell $ws.Range("D22") '232.85'
Set-TextCell $ws.Range("E22") '  +2.32%  '
Set-TextCell $ws.Range("E23") '  +0.00%  '
Set-TextCell $ws.Range("E24") '  -0.49%  '
Set-TextCell $ws.Range("D25") '2.38'
Set-TextCell $ws.Range("E25") '  +2.18%  '
Set-TextCell $ws.Range("D26") '9.76'
Set-TextCell $ws.Range("E26") '  +1.85%  '
Set-TextCell $ws.Range("D27") '172.64'
Set-TextCell $ws.Range("E27") '  +0.96%  '
Set-TextCell $ws.Range("E28") '  -0.39%  '
Set-TextCell $ws.Range("D29") '20.14'
Set-TextCell $ws.Range("E29") '  +4.59%  '
Set-TextCell $ws.Range("D30") '1.41'
Set-TextCell $ws.Range("E30") '  -1.98%  '
Set-TextCell $ws.Range("D31") '2.77'
Set-TextCell $ws.Range("E31") '  +12.35%  '
Set-TextCell $ws.Range("E32") '  +1.77%  '
Set-TextCell $ws.Range("D33") '4.67'
Set-TextCell $ws.Range("E33") '  +3.09%  '
Set-TextCell $ws.Range("D34") '4.84'
Set-TextCell $ws.Range("E34") '  +2.61%  '
Set-TextCell $ws.Range("D35") '7.20'
Set-TextCell $ws.Range("E35") '  +10.85%  '
Set-TextCell $ws.Range("D36") '0.0624'
Set-TextCell $ws.Range("E36") '  +1.99%  '
Set-TextCell $ws.Range("D37") '2.45'
Set-TextCell $ws.Range("E37") '  +2.17%  '
Set-TextCell $ws.Range("E38") '  +0.55%  '
Set-TextCell $ws.Range("D39") '1.00'
Set-TextCell $ws.Range("E39") '  +0.17%  '
Set-TextCell $ws.Range("D40") '105.12'
Set-TextCell $ws.Range("E40") '  +3.82%  '
Set-TextCell $ws.Range("E41") '  +0.39%  '
Set-TextCell $ws.Range("D42") '18.05'
Set-TextCell $ws.Range("E42") '  -0.19%  '
Set-TextCell $ws.Range("D43") '1.541.38'
Set-TextCell $ws.Range("E43") '  +0.25%  '
Set-TextCell $ws.Range("E44") '  +5.88%  '
Set-TextCell $ws.Range("B45") 'ARBITRUM'
Set-TextCell $ws.Range("C45") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range("D45") '1.12'
Set-TextCell $ws.Range("E45") '  +8.08%  '
Set-TextCell $ws.Range("B46") 'FTXToken'
Set-TextCell $ws.Range("C46") 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell $ws.Range("D46") '4.31'
Set-TextCell $ws.Range("E46") '  +5.04%  '
Set-TextCell $ws.Range("B47") 'Cronos'
Set-TextCell $ws.Range("C47") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range("D47") '0.0929'
Set-TextCell $ws.Range("E47") '  +1.19%  '
Set-TextCell $ws.Range("B48") 'HuobiToken'
Set-TextCell $ws.Range("C48") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws.Range("D48") '2.83'
Set-TextCell $ws.Range("E48") '  +0.36%  '
Set-TextCell $ws.Range("B49") 'FraxShare'
Set-TextCell $ws.Range("C49") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range("D49") '7.85'
Set-TextCell $ws.Range("E49") '  +1.12%  '
Set-TextCell $ws.Range("D50") '2.381.10'
Set-TextCell $ws.Range("E50") '  +3.96%  '
Set-TextCell $ws.Range("E51") '  +0.25%  '
